$wb = $excel.ActiveWorkbook

# 1. Rename sheet "buscar" -> "buscarHome"
$wsBuscar = $wb.Worksheets.Item("buscar")
$wsBuscar.Name = "buscarHome"

# 2. Update the product list on the buscarHome sheet
#    Row2 gets a new product (with word-wrap turned on), rows 3 & 4 are cleared out.
$wsBuscar.Range("A2").Value = "BOSE SOUNDLINK BLUETOOTH SPEAKER III"
$wsBuscar.Range("A2").WrapText = $true
$wsBuscar.Range("A3").ClearContents()
$wsBuscar.Range("A4").ClearContents()

# 3. Select A2 on buscarHome and make this the active/selected sheet/tab
$wsBuscar.Range("A2").Select()
$wsBuscar.Activate()

$wb.Save()
